$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.03958866666666667
$ws.Range("H2").Value = 0.118766
$ws.Range("I2").Value = 0.0007442768123675562
$ws.Range("J2").Value = 0.0007442768123675561
$ws.Range("M2").Value = 0.1419263333333333
$ws.Range("N2").Value = 0.425779
$ws.Range("O2").Value = 0.002583058778296354
$ws.Range("P2").Value = 0.002583058778296354
$ws.Range("Q2").Value = 0.005618674301555556
$ws.Range("R2").Value = 0.050568068714
$ws.Range("S2").Value = 0.000001922510753668445
$ws.Range("T2").Value = 0.000001922510753668444

$ws.Range("G3").Value = 0.03958866666666667
$ws.Range("H3").Value = 0.118766
$ws.Range("I3").Value = 0.0007442768123675562
$ws.Range("J3").Value = 0.0007442768123675561
$ws.Range("O3").Value = 0.001399682868699959
$ws.Range("P3").Value = 0.001399682868699959
$ws.Range("Q3").Value = 0.003044592802444444
$ws.Range("R3").Value = 0.027401335222
$ws.Range("S3").Value = 0.000001041751503841482
$ws.Range("T3").Value = 0.000001041751503841482

$ws.Range("G4").Value = 0.03958866666666667
$ws.Range("H4").Value = 0.118766
$ws.Range("I4").Value = 0.0007442768123675562
$ws.Range("J4").Value = 0.0007442768123675561
$ws.Range("M4").Value = 2.613991
$ws.Range("N4").Value = 7.841973
$ws.Range("O4").Value = 0.04757462720522382
$ws.Range("P4").Value = 0.04757462720522382
$ws.Range("Q4").Value = 0.1034844183686667
$ws.Range("R4").Value = 0.931359765318
$ws.Range("S4").Value = 0.00003540869188587881
$ws.Range("T4").Value = 0.0000354086918858788

$ws.Range("G5").Value = 0.03958866666666667
$ws.Range("H5").Value = 0.118766
$ws.Range("I5").Value = 0.0007442768123675562
$ws.Range("J5").Value = 0.0007442768123675561
$ws.Range("M5").Value = 52.11224233333333
$ws.Range("N5").Value = 156.336727
$ws.Range("O5").Value = 0.9484426311477799
$ws.Range("P5").Value = 0.9484426311477798
$ws.Range("Q5").Value = 2.063054190986889
$ws.Range("R5").Value = 18.567487718882
$ws.Range("S5").Value = 0.0007059038582241675
$ws.Range("T5").Value = 0.0007059038582241674

$ws.Range("H6").Value = 0.059669
$ws.Range("I6").Value = 0.0003739306966401135
$ws.Range("J6").Value = 0.0003739306966401134
$ws.Range("M6").Value = 0.1419263333333333
$ws.Range("N6").Value = 0.425779
$ws.Range("O6").Value = 0.002583058778296354
$ws.Range("P6").Value = 0.002583058778296354
$ws.Range("Q6").Value = 0.002822867461222222
$ws.Range("R6").Value = 0.025405807151
$ws.Range("S6").Value = 0.0000009658849684307161
$ws.Range("T6").Value = 0.0000009658849684307157

$ws.Range("H7").Value = 0.059669
$ws.Range("I7").Value = 0.0003739306966401135
$ws.Range("J7").Value = 0.0003739306966401134
$ws.Range("O7").Value = 0.001399682868699959
$ws.Range("P7").Value = 0.001399682868699959
$ws.Range("S7").Value = 0.0000005233843901682081
$ws.Range("T7").Value = 0.000000523384390168208

$ws.Range("H8").Value = 0.059669
$ws.Range("I8").Value = 0.0003739306966401135
$ws.Range("J8").Value = 0.0003739306966401134
$ws.Range("M8").Value = 2.613991
$ws.Range("N8").Value = 7.841973
$ws.Range("O8").Value = 0.04757462720522382
$ws.Range("P8").Value = 0.04757462720522382
$ws.Range("Q8").Value = 0.05199140965966666
$ws.Range("R8").Value = 0.467922686937
$ws.Range("S8").Value = 0.00001778961349324304
$ws.Range("T8").Value = 0.00001778961349324303

$ws.Range("H9").Value = 0.059669
$ws.Range("I9").Value = 0.0003739306966401135
$ws.Range("J9").Value = 0.0003739306966401134
$ws.Range("M9").Value = 52.11224233333333
$ws.Range("N9").Value = 156.336727
$ws.Range("O9").Value = 0.9484426311477799
$ws.Range("P9").Value = 0.9484426311477798
$ws.Range("Q9").Value = 1.036495129262555
$ws.Range("R9").Value = 9.328456163363001
$ws.Range("S9").Value = 0.0003546518137882715
$ws.Range("T9").Value = 0.0003546518137882714

$ws.Range("G10").Value = 53.131305
$ws.Range("H10").Value = 159.393915
$ws.Range("I10").Value = 0.9988817924909924
$ws.Range("J10").Value = 0.9988817924909923
$ws.Range("M10").Value = 0.1419263333333333
$ws.Range("N10").Value = 0.425779
$ws.Range("O10").Value = 0.002583058778296354
$ws.Range("P10").Value = 0.002583058778296354
$ws.Range("Q10").Value = 7.540731303865001
$ws.Range("R10").Value = 67.866581734785
$ws.Range("S10").Value = 0.002580170382574255
$ws.Range("T10").Value = 0.002580170382574255

$ws.Range("G11").Value = 53.131305
$ws.Range("H11").Value = 159.393915
$ws.Range("I11").Value = 0.9988817924909924
$ws.Range("J11").Value = 0.9988817924909923
$ws.Range("O11").Value = 0.001399682868699959
$ws.Range("P11").Value = 0.001399682868699959
$ws.Range("Q11").Value = 4.086098431895
$ws.Range("R11").Value = 36.774885887055
$ws.Range("S11").Value = 0.001398117732805949
$ws.Range("T11").Value = 0.001398117732805949

$ws.Range("G12").Value = 53.131305
$ws.Range("H12").Value = 159.393915
$ws.Range("I12").Value = 0.9988817924909924
$ws.Range("J12").Value = 0.9988817924909923
$ws.Range("M12").Value = 2.613991
$ws.Range("N12").Value = 7.841973
$ws.Range("O12").Value = 0.04757462720522382
$ws.Range("P12").Value = 0.04757462720522382
$ws.Range("Q12").Value = 138.884753088255
$ws.Range("R12").Value = 1249.962777794295
$ws.Range("S12").Value = 0.0475214288998447
$ws.Range("T12").Value = 0.0475214288998447

$ws.Range("G13").Value = 53.131305
$ws.Range("H13").Value = 159.393915
$ws.Range("I13").Value = 0.9988817924909924
$ws.Range("J13").Value = 0.9988817924909923
$ws.Range("M13").Value = 52.11224233333333
$ws.Range("N13").Value = 156.336727
$ws.Range("O13").Value = 0.9484426311477799
$ws.Range("P13").Value = 0.9484426311477798
$ws.Range("Q13").Value = 2768.791441646245
$ws.Range("R13").Value = 24919.1229748162
$ws.Range("S13").Value = 0.9473820754757676
$ws.Range("T13").Value = 0.9473820754757674
